# Auto-generated edit script applying the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.839.17"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.543.34"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'205.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "'21.38"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Value = "'0.0852"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.30%  "
$ws.Range("D12").Value = "1.762.80"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Value = "1.547.08"
$ws.Range("E13").Value = "  -1.25%  "
$ws.Range("E14").Value = "  -1.49%  "
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("D16").Value = "26.840.57"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'61.30"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").Value = "'215.21"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0682"
$ws.Range("E19").Value = "  +0.57%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.22"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'4.01"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("D23").Value = "'9.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.82%  "
$ws.Range("E24").Value = "  -2.80%  "
$ws.Range("D25").Value = "'152.86"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").Value = "'6.61"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").Value = "'14.84"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "'0.0458"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.10"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.59%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "1.365.40"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "'0.958"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.31%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "'0.0164"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "'0.806"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "'5.69"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.89%  "
$ws.Range("D43").Value = "'0.989"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("D45").Value = "'63.13"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  -3.85%  "
$ws.Range("D47").Value = "1.677.33"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").Value = "'84.21"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "'0.0513"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.31%  "
$ws.Range("D50").Value = "0.0₇0965"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("E51").Value = "  +0.14%  "
